{"js": "// Update the \"Journal article\" row's UK count (764 -> 763) and its row\n// total (1,471 -> 1,470), plus the corresponding \"Total\" row's UK column\n// (1,146 -> 1,145) and grand total (2,308 -> 2,307) in the single summary\n// table in the document body. Each old value is a unique, whole table-cell\n// text, so an exact (case-sensitive, whole-match) search-and-replace on the\n// body is unambiguous.\n\nconst replacements = [\n  [\"764\", \"763\"],\n  [\"1,471\", \"1,470\"],\n  [\"1,146\", \"1,145\"],\n  [\"2,308\", \"2,307\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${oldText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the \"Journal article\" row's UK count (764 -> 763) and its row\n# total (1,471 -> 1,470), plus the corresponding \"Total\" row's UK column\n# (1,146 -> 1,145) and grand total (2,308 -> 2,307) in the single summary\n# table in the document body. Each old value is unique in the document, so\n# an exact (case-sensitive, whole-word) Find/Replace is unambiguous.\n\n$d = $word.ActiveDocument\n\nfunction Replace-UniqueText($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $findText\n    $range.Find.Replacement.Text = $replaceText\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.Format = $false\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $true\n    $range.Find.MatchWildcards = $false\n\n    $ok = $range.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Find/Replace failed: could not find '$findText'\"\n    }\n}\n\nReplace-UniqueText \"764\" \"763\"\nReplace-UniqueText \"1,471\" \"1,470\"\nReplace-UniqueText \"1,146\" \"1,145\"\nReplace-UniqueText \"2,308\" \"2,307\"\n"}
